$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.913.19"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.660.55"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -1.09%  "
$ws.Range("D5").Value = "'324.61"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.3640"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'47.48"
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "'0.3271"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "'1.136"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'0.07090"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'0.9989"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'6.063"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "'19.60"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "1.658.50"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "'6.611"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "'0.00001050"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "'0.06615"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'0.9988"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").Value = "'79.19"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").Value = "'5.927"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "'15.81"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "24.884.48"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "'2.438"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'2.449"
$ws.Range("E26").Value = "  -6.38%  "
$ws.Range("D27").Value = "'147.76"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'18.68"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("D29").Value = "1.840.50"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'1.205"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'125.38"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'4.101"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "'5.780"
$ws.Range("E33").Value = "  -7.67%  "
$ws.Range("D34").Value = "'0.08480"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "'1.639"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Value = "'12.30"
$ws.Range("E36").Value = "  -5.83%  "
$ws.Range("D37").Value = "'1.279"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").Value = "'5.180"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'0.02272"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("D40").Value = "'0.06080"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").Value = "'8.387"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").Value = "'0.2073"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'0.9991"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "'0.5955"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").Value = "'13.84"
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("D46").Value = "'3.865"
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("D47").Value = "'0.5635"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'125.14"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'1.953"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").Value = "'0.06984"
$ws.Range("E50").Value = "  -1.35%  "
$ws.Range("D51").Value = "'1.191"
$ws.Range("E51").Value = "  +0.75%  "
